# validacion de rubros
# Adds new "i_errormsg"/"i_isimported"/"processing"/"processed" parameter rows
# to Hoja3, makes Hoja3 the active/selected sheet (instead of Hoja4).

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Hoja3 becomes the selected/active sheet (was Hoja4 before).
$ws3.Activate()

# New validation/rubro parameter rows.
$ws3.Range("A4").Value = "i_errormsg"
$ws3.Range("B4").Value = "CV "
$ws3.Range("C4").Value = 2000

$ws3.Range("A5").Value = "i_isimported"
$ws3.Range("B5").Value = "CHARACTER"
$ws3.Range("C5").Value = 1

$ws3.Range("A6").Value = "processing"
$ws3.Range("B6").Value = "CHARACTER"
$ws3.Range("C6").Value = 1

$ws3.Range("A7").Value = "processed"
$ws3.Range("B7").Value = "CHARACTER"
$ws3.Range("C7").Value = 1
$ws3.Range("D7").Value = "N"

# Column A is best-fit/widened to accommodate the new content.
$ws3.Columns.Item(1).ColumnWidth = 14.3

# Final selection sits on the last new row, matching the edited workbook.
$ws3.Range("A7").Select()
